# ISES2023 poster: ensure consistency of SSA terminology
# ("Non-Targeted Analysis" / "NTA" -> "Suspect Screening Analysis" / "SSA")
#
# Touches three shapes on slide 1:
#   id=8  "TextBox 7"  - poster title
#   id=22 "TextBox 21" - Methods bullet
#   id=40 "TextBox 39" - Figure 1 caption

function Get-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title (shape id 8): shrink font and retitle ---
$titleShape = Get-ShapeById $s.Shapes 8
$trTitle = $titleShape.TextFrame.TextRange
# Set the new font size *before* re-writing the text so the text-triggered
# autofit relayout (spAutoFit) measures the new title at its final size.
$trTitle.Font.Size = 78
$trTitle.Text = "Developing Chemical Signatures for 5 Categories of Household Products Using Suspect Screening Analysis"

# --- Methods bullet (shape id 22): "Non-targeted" -> "Suspect screening" ---
$methodsShape = Get-ShapeById $s.Shapes 22
$trMethods = $methodsShape.TextFrame.TextRange
$methodsPara = $trMethods.Paragraphs(2)
$methodsRun = $methodsPara.Runs(1)
# First 12 characters of the run are exactly "Non-targeted"; replace them.
$methodsPrefix = $methodsRun.Characters(1, 12)
$methodsPrefix.Text = "Suspect screening"

# --- Figure 1 caption (shape id 40): "non-targeted" -> "suspect screening" ---
$figShape = Get-ShapeById $s.Shapes 40
$trFig = $figShape.TextFrame.TextRange
$figPara = $trFig.Paragraphs(1)
$figRun = $figPara.Runs(2)
$figRun.Text = "Workflow of suspect screening analysis of products from 5 types of household consumer products. Products were extracted with dichloromethane (DCM). After addition of an internal standard, each extraction was analyzed via GC X GC-TOFMS to obtain its mass spectra. The spectra were matched to the 2017 NIST database and analytical standards were used to confirm a subset of the chemical identifications. Chemicals were annotated by reported or predicted functional uses"
